$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 8 (shifts existing rows 8-40 down to 9-41)
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new weekly price record
$ws.Range("A8").Value = 8
$ws.Range("B8").Value = "Terminal La Palmera de La Serena"
$ws.Range("C8").Value = "Coquimbo"
$ws.Range("D8").Value = 44701
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 100114007
$ws.Range("G8").Value = "Jengibre"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 440
$ws.Range("K8").Value = 14000
$ws.Range("L8").Value = 15000
$ws.Range("M8").Value = 14500
$ws.Range("N8").Value = "`$/caja 13 kilos"
$ws.Range("O8").Value = "Perú"
$ws.Range("P8").Value = 1115
$ws.Range("Q8").Value = 13
$ws.Range("R8").Value = "Hortaliza"
